# Applies the "Added a few more slots" edit to the Arcane Gems review:
#   1. Insert a new "Meta description: ..." paragraph right after the
#      Heading1 title paragraph ("Play Arcane Gems for Free- Review of
#      this High-Rewarding Slot").
#   2. At the very end of the document, drop the paragraph that duplicates
#      that same bold title text, and replace the text of the remaining
#      (italic) paragraph with the image-generation "Prompt: ..." text,
#      keeping its italic run formatting intact.

$d = $word.ActiveDocument

$titleText    = "Play Arcane Gems for Free- Review of this High-Rewarding Slot"
$metaBoldText = "Meta description"
$metaRestText = ": Read our review of Arcane Gems and play for free. Find out about its high rewards, symbol locking feature, and respin feature. Available for free play."
$oldTailText  = "Read our review of Arcane Gems and play for free. Find out about its high rewards, symbol locking feature, and respin feature. Available for free play."
$promptText   = 'Prompt: Create a cartoon-style feature image for the game "Arcane Gems" that features a happy Maya warrior with glasses. For this feature image, we want to bring in elements of both the theme of gems and the unique aspect of the game''s respin feature. The Maya warrior with glasses will add a touch of personality to the image and make it stand out. The Maya warrior should be depicted with a big smile on their face, eyes twinkling behind their glasses. They should be surrounded by piles of colorful gems, with one hand clutching a handful of gems, and the other hand pointing to the reels of the game. The reels should be shown on the image, with the game name "Arcane Gems" prominently displayed. The symbols on the reels should be vibrant and eye-catching, particularly highlighting the blue gem symbol that pays out the most. The symbol locking feature should also be represented, perhaps with the locked symbols being depicted as glowing and surrounded by a blue aura. Overall, the image should be lively and fun, drawing players in with its bright colors, charming character, and attention to the unique features of the game.'

# ---------------------------------------------------------------------
# Locate the first paragraph (the Heading1 title) and the two paragraphs
# at the end of the body (the duplicated bold title + the italic blurb).
# ---------------------------------------------------------------------
$titleParaIndex = 0
$boldDupIndex = 0
$italicParaIndex = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($titleParaIndex -eq 0 -and $t.TrimEnd() -eq $titleText) {
        $titleParaIndex = $i
    }
    if ($t.TrimEnd() -eq $titleText -and $i -ne $titleParaIndex) {
        $boldDupIndex = $i
    }
    if ($t.TrimEnd() -eq $oldTailText) {
        $italicParaIndex = $i
    }
}

# ---------------------------------------------------------------------
# 1) Insert the "Meta description" paragraph right after the title.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item($titleParaIndex)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item($titleParaIndex + 1)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.InsertBefore($metaBoldText + $metaRestText)

$metaPara = $d.Paragraphs.Item($titleParaIndex + 1)
$metaStart = $metaPara.Range.Start
$boldRange = $d.Range($metaStart, $metaStart + $metaBoldText.Length)
$boldRange.Bold = 1

# Inserting the paragraph shifted every later paragraph's index by 1.
$boldDupIndex = $boldDupIndex + 1
$italicParaIndex = $italicParaIndex + 1

# ---------------------------------------------------------------------
# 2) Delete the duplicated bold title paragraph, then replace the
#    italic paragraph's text with the image prompt (keeping its
#    existing italic run formatting).
# ---------------------------------------------------------------------
$boldTitlePara = $d.Paragraphs.Item($boldDupIndex)
$boldTitlePara.Range.Delete()

if ($boldDupIndex -lt $italicParaIndex) {
    $italicParaIndex = $italicParaIndex - 1
}

$italicPara = $d.Paragraphs.Item($italicParaIndex)
$italicRange = $italicPara.Range
$textOnlyRange = $d.Range($italicRange.Start, $italicRange.End - 1)
$textOnlyRange.Text = $promptText

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
